$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "Meta description: ..." paragraph that sits right after the
#    "Play Cash Coaster slot game for free" Heading1 paragraph.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Meta description")) {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 2. Insert a new bold paragraph "Play Cash Coaster slot game for free"
#    right before the last paragraph of the document (the one that used to
#    hold the "Create a feature image ..." image prompt).
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$beforeLast = $d.Paragraphs.Item($count - 1)
# Split a couple of characters into the previous paragraph's range so the
# inserted OOXML paragraph lands cleanly *between* the two existing
# paragraphs instead of merging with / replacing either of them.
$beforeLastLen = $beforeLast.Range.End - $beforeLast.Range.Start
$offset = [Math]::Min(3, [Math]::Max(1, $beforeLastLen - 1))
$splitPos = $beforeLast.Range.Start + $offset
$insertionPoint = $d.Range($splitPos, $splitPos)
$insertionPoint.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Cash Coaster slot game for free</w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# 3. Replace the text of the (now) final paragraph - the old image-prompt
#    text - with the meta-description copy, keeping its italic formatting.
# ---------------------------------------------------------------------------
$oldText = "Create a feature image for Cash Coaster that captures the fun, upbeat roller coaster theme of the game. The image should be in cartoon style and feature a happy Maya warrior with glasses, who represents the excitement and thrill of the amusement park. This warrior should be shown riding a roller coaster with a big smile on their face, while holding some of the classic amusement park treats like pretzels, cotton candy, and caramel apples. The background could include the roller coaster and the bright neon lights of the Cash Coaster logo, as well as other carnival attractions like a Ferris wheel or a carousel. Overall, the image should convey the playful and exciting vibe of Cash Coaster and entice players to take a ride on this thrilling slot game."
$newText = "Read our review of Cash Coaster, a 5-reel, 30-payline slot game with Wild, Scatter, and Free Spins features. Play for free and potentially win big."

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastPara.Range
$lastRange.Find.ClearFormatting()
$lastRange.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

Write-Output "done"
